# Auto-generated: update market-price derived columns (H-N) per scheduled runner refresh
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7990.731
$ws.Range("I40").Value = 5810.6
$ws.Range("J40").Value = 9353.3125
$ws.Range("K40").Value = 5810.6
$ws.Range("L40").Value = 9353.3125
$ws.Range("M40").Value = -5635.6
$ws.Range("N40").Value = -9703.3125
$ws.Range("H45").Value = 3902.8333
$ws.Range("J45").Value = 8607.4
$ws.Range("L45").Value = 25822.2
$ws.Range("N45").Value = -26206.2
$ws.Range("H80").Value = 772.9167
$ws.Range("I80").Value = 655.4
$ws.Range("J80").Value = 856.8570999999999
$ws.Range("K80").Value = 1966.2
$ws.Range("L80").Value = 2570.5713
$ws.Range("M80").Value = -968.1999999999998
$ws.Range("N80").Value = -4566.5713
$ws.Range("H83").Value = 772.9167
$ws.Range("I83").Value = 655.4
$ws.Range("J83").Value = 856.8570999999999
$ws.Range("K83").Value = 5898.599999999999
$ws.Range("L83").Value = 7711.7139
$ws.Range("M83").Value = -906.5999999999995
$ws.Range("N83").Value = -17695.7139
$ws.Range("H138").Value = 4475.7295
$ws.Range("I138").Value = 2726.7778
$ws.Range("J138").Value = 6132.6313
$ws.Range("K138").Value = 8180.3334
$ws.Range("L138").Value = 18397.8939
$ws.Range("M138").Value = -3040.3334
$ws.Range("N138").Value = -28677.8939
$ws.Range("H141").Value = 3874.6191
$ws.Range("I141").Value = 2347.6155
$ws.Range("K141").Value = 7042.8465
$ws.Range("M141").Value = -1862.8465

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 4000.8572
$ws.Range("I26").Value = 1601.2
$ws.Range("K26").Value = 1601.2
$ws.Range("M26").Value = -1271.2
$ws.Range("H32").Value = 1514.2063
$ws.Range("I32").Value = 1368.7455
$ws.Range("K32").Value = 1368.7455
$ws.Range("M32").Value = -1081.7455
$ws.Range("H39").Value = 2142.8572
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 5000
$ws.Range("N39").Value = -6040
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("H61").Value = 3339.2632
$ws.Range("I61").Value = 2462.6365
$ws.Range("K61").Value = 2462.6365
$ws.Range("M61").Value = -2250.6365
$ws.Range("H74").Value = 8132360
$ws.Range("I74").Value = 9261459
$ws.Range("J74").Value = 2846
$ws.Range("K74").Value = 9261459
$ws.Range("L74").Value = 2846
$ws.Range("M74").Value = -9260585
$ws.Range("N74").Value = -4594
$ws.Range("H77").Value = 8132360
$ws.Range("I77").Value = 9261459
$ws.Range("J77").Value = 2846
$ws.Range("K77").Value = 46307295
$ws.Range("L77").Value = 14230
$ws.Range("M77").Value = -46302927
$ws.Range("N77").Value = -22966
$ws.Range("H102").Value = 1516.5
$ws.Range("I102").Value = 1516.5
$ws.Range("K102").Value = 1516.5
$ws.Range("M102").Value = 105.5
$ws.Range("H122").Value = 3964.9443
$ws.Range("I122").Value = 2937
$ws.Range("J122").Value = 5249.875
$ws.Range("K122").Value = 8811
$ws.Range("L122").Value = 15749.625
$ws.Range("M122").Value = -6361
$ws.Range("N122").Value = -20649.625
$ws.Range("H132").Value = 1967.7368
$ws.Range("I132").Value = 974.8823
$ws.Range("K132").Value = 2924.6469
$ws.Range("M132").Value = -394.6468999999997
$ws.Range("H136").Value = 3339.2632
$ws.Range("I136").Value = 2462.6365
$ws.Range("K136").Value = 7387.9095
$ws.Range("M136").Value = -4837.9095

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2319.2
$ws.Range("I107").Value = 1400.25
$ws.Range("K107").Value = 1400.25
$ws.Range("M107").Value = 519.75
$ws.Range("H134").Value = 5289.16
$ws.Range("I134").Value = 3478.9443
$ws.Range("K134").Value = 10436.8329
$ws.Range("M134").Value = -7901.832900000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1000072.75
$ws.Range("I19").Value = 1176491.5
$ws.Range("K19").Value = 1176491.5
$ws.Range("M19").Value = -1176321.5
$ws.Range("H24").Value = 1000072.75
$ws.Range("I24").Value = 1176491.5
$ws.Range("K24").Value = 1176491.5
$ws.Range("M24").Value = -1176321.5
$ws.Range("H28").Value = 44833.332
$ws.Range("J28").Value = 44833.332
$ws.Range("L28").Value = 44833.332
$ws.Range("N28").Value = -45323.332
$ws.Range("H58").Value = 4468.04
$ws.Range("I58").Value = 3300.718
$ws.Range("K58").Value = 3300.718
$ws.Range("M58").Value = -3097.718
$ws.Range("H86").Value = 4426.516
$ws.Range("I86").Value = 2579.5
$ws.Range("K86").Value = 2579.5
$ws.Range("M86").Value = -1456.5
$ws.Range("H89").Value = 4426.516
$ws.Range("I89").Value = 2579.5
$ws.Range("K89").Value = 12897.5
$ws.Range("M89").Value = -7281.5
$ws.Range("H93").Value = 16369.25
$ws.Range("I93").Value = 17412.334
$ws.Range("K93").Value = 17412.334
$ws.Range("M93").Value = -15540.334
$ws.Range("H95").Value = 24038.166
$ws.Range("J95").Value = 24038.166
$ws.Range("L95").Value = 24038.166
$ws.Range("N95").Value = -29530.166
$ws.Range("H96").Value = 11500
$ws.Range("J96").Value = 11500
$ws.Range("L96").Value = 11500
$ws.Range("N96").Value = -16992
$ws.Range("H132").Value = 3609.8276
$ws.Range("I132").Value = 2827.24
$ws.Range("J132").Value = 8501
$ws.Range("K132").Value = 8481.719999999999
$ws.Range("L132").Value = 25503
$ws.Range("M132").Value = -5951.719999999999
$ws.Range("N132").Value = -30563
$ws.Range("H134").Value = 4739.1333
$ws.Range("I134").Value = 2260.25
$ws.Range("K134").Value = 6780.75
$ws.Range("M134").Value = -4245.75
$ws.Range("H136").Value = 4468.04
$ws.Range("I136").Value = 3300.718
$ws.Range("K136").Value = 9902.153999999999
$ws.Range("M136").Value = -7352.153999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1721.8334
$ws.Range("I136").Value = 1373.3636
$ws.Range("K136").Value = 4120.0908
$ws.Range("M136").Value = 979.9092000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 31335.5
$ws.Range("J106").Value = 31335.5
$ws.Range("L106").Value = 31335.5
$ws.Range("N106").Value = -33859.5
$ws.Range("H132").Value = 4318.722
$ws.Range("I132").Value = 1955.3478
$ws.Range("J132").Value = 8500.076999999999
$ws.Range("K132").Value = 5866.0434
$ws.Range("L132").Value = 25500.231
$ws.Range("M132").Value = -3336.0434
$ws.Range("N132").Value = -30560.231
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = ""
$ws.Range("N138").Value = 0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2923.111
$ws.Range("I16").Value = 2724.3635
$ws.Range("K16").Value = 2724.3635
$ws.Range("M16").Value = -2554.3635
$ws.Range("H22").Value = 2227.3333
$ws.Range("I22").Value = 1279.4
$ws.Range("J22").Value = 3412.25
$ws.Range("K22").Value = 1279.4
$ws.Range("L22").Value = 3412.25
$ws.Range("M22").Value = -984.4000000000001
$ws.Range("N22").Value = -4002.25
$ws.Range("H27").Value = 2227.3333
$ws.Range("I27").Value = 1279.4
$ws.Range("J27").Value = 3412.25
$ws.Range("K27").Value = 1279.4
$ws.Range("L27").Value = 3412.25
$ws.Range("M27").Value = -1172.4
$ws.Range("N27").Value = -3626.25
$ws.Range("H38").Value = 100000
$ws.Range("J38").Value = 100000
$ws.Range("L38").Value = 100000
$ws.Range("N38").Value = -100820
$ws.Range("H46").Value = 2809
$ws.Range("I46").Value = 1160
$ws.Range("J46").Value = 3908.3333
$ws.Range("K46").Value = 1160
$ws.Range("L46").Value = 3908.3333
$ws.Range("M46").Value = -972
$ws.Range("N46").Value = -4284.3333
$ws.Range("H55").Value = 5555949.5
$ws.Range("I55").Value = 10000319
$ws.Range("K55").Value = 10000319
$ws.Range("M55").Value = -10000146
$ws.Range("H93").Value = 15681.342
$ws.Range("I93").Value = 12202.667
$ws.Range("K93").Value = 12202.667
$ws.Range("M93").Value = -10954.667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8334
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 11501
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 23002
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -25124
$ws.Range("H84").Value = 8334
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 11501
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 115010
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -125618
$ws.Range("H132").Value = 4847.2876
$ws.Range("I132").Value = 3827.1929
$ws.Range("J132").Value = 8481.375
$ws.Range("K132").Value = 11481.5787
$ws.Range("L132").Value = 25444.125
$ws.Range("M132").Value = -8951.5787
$ws.Range("N132").Value = -30504.125
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
$ws.Range("H137").Value = 63386
$ws.Range("J137").Value = 63386
$ws.Range("L137").Value = 63386
$ws.Range("N137").Value = -73586
$ws.Range("H141").Value = 141618.33
$ws.Range("J141").Value = 141618.33
$ws.Range("L141").Value = 141618.33
$ws.Range("N141").Value = -151978.33

Write-Host "Applied market-data refresh across 8 sheets"